{"js": "// Update the 25 two-digit multiplication prompts in the practice table.\n// Each old expression is unique in the document, so searching the body\n// for each one and replacing its single match is safe and order-independent.\n\nconst pairs = [\n  [\"72\u00d783=\", \"33\u00d730=\"],\n  [\"99\u00d752=\", \"59\u00d739=\"],\n  [\"35\u00d717=\", \"73\u00d737=\"],\n  [\"16\u00d772=\", \"89\u00d734=\"],\n  [\"20\u00d721=\", \"99\u00d790=\"],\n  [\"59\u00d782=\", \"79\u00d732=\"],\n  [\"30\u00d735=\", \"97\u00d792=\"],\n  [\"89\u00d729=\", \"92\u00d740=\"],\n  [\"28\u00d717=\", \"82\u00d767=\"],\n  [\"87\u00d746=\", \"98\u00d739=\"],\n  [\"96\u00d763=\", \"43\u00d765=\"],\n  [\"96\u00d746=\", \"89\u00d762=\"],\n  [\"67\u00d759=\", \"53\u00d784=\"],\n  [\"52\u00d778=\", \"17\u00d731=\"],\n  [\"43\u00d717=\", \"71\u00d793=\"],\n  [\"23\u00d776=\", \"29\u00d787=\"],\n  [\"26\u00d712=\", \"19\u00d760=\"],\n  [\"53\u00d761=\", \"99\u00d765=\"],\n  [\"99\u00d778=\", \"17\u00d769=\"],\n  [\"37\u00d795=\", \"30\u00d718=\"],\n  [\"13\u00d745=\", \"90\u00d715=\"],\n  [\"24\u00d729=\", \"70\u00d762=\"],\n  [\"39\u00d749=\", \"58\u00d781=\"],\n  [\"96\u00d783=\", \"89\u00d749=\"],\n  [\"44\u00d723=\", \"97\u00d781=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 two-digit multiplication prompts in the practice table.\n# Each old expression is unique in the document, so a simple Find/Replace\n# (wdReplaceOne) per pair is safe and order-independent.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"72\u00d783=\", \"33\u00d730=\"),\n    @(\"99\u00d752=\", \"59\u00d739=\"),\n    @(\"35\u00d717=\", \"73\u00d737=\"),\n    @(\"16\u00d772=\", \"89\u00d734=\"),\n    @(\"20\u00d721=\", \"99\u00d790=\"),\n    @(\"59\u00d782=\", \"79\u00d732=\"),\n    @(\"30\u00d735=\", \"97\u00d792=\"),\n    @(\"89\u00d729=\", \"92\u00d740=\"),\n    @(\"28\u00d717=\", \"82\u00d767=\"),\n    @(\"87\u00d746=\", \"98\u00d739=\"),\n    @(\"96\u00d763=\", \"43\u00d765=\"),\n    @(\"96\u00d746=\", \"89\u00d762=\"),\n    @(\"67\u00d759=\", \"53\u00d784=\"),\n    @(\"52\u00d778=\", \"17\u00d731=\"),\n    @(\"43\u00d717=\", \"71\u00d793=\"),\n    @(\"23\u00d776=\", \"29\u00d787=\"),\n    @(\"26\u00d712=\", \"19\u00d760=\"),\n    @(\"53\u00d761=\", \"99\u00d765=\"),\n    @(\"99\u00d778=\", \"17\u00d769=\"),\n    @(\"37\u00d795=\", \"30\u00d718=\"),\n    @(\"13\u00d745=\", \"90\u00d715=\"),\n    @(\"24\u00d729=\", \"70\u00d762=\"),\n    @(\"39\u00d749=\", \"58\u00d781=\"),\n    @(\"96\u00d783=\", \"89\u00d749=\"),\n    @(\"44\u00d723=\", \"97\u00d781=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
